$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 504.1111
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()
$ws.Range("H40").Value = 2229.1177
$ws.Range("I40").Value = 2049.7856
$ws.Range("J40").Value = 3066
$ws.Range("K40").Value = 2049.7856
$ws.Range("L40").Value = 3066
$ws.Range("M40").Value = -1874.7856
$ws.Range("N40").Value = -3416
$ws.Range("H51").Value = 259999.5
$ws.Range("J51").Value = 259999.5
$ws.Range("L51").Value = 259999.5
$ws.Range("N51").Value = -260967.5
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("N67").ClearContents()
$ws.Range("H70").Value = 3008.6956
$ws.Range("I70").Value = 3104.762
$ws.Range("K70").Value = 9314.286
$ws.Range("M70").Value = -9044.286
$ws.Range("H73").Value = 3008.6956
$ws.Range("I73").Value = 3104.762
$ws.Range("K73").Value = 9314.286
$ws.Range("M73").Value = -8378.286
$ws.Range("H74").Value = 2000
$ws.Range("I74").Value = 2000
$ws.Range("K74").Value = 2000
$ws.Range("M74").Value = -1064
$ws.Range("H77").Value = 2000
$ws.Range("I77").Value = 2000
$ws.Range("K77").Value = 10000
$ws.Range("M77").Value = -5320
$ws.Range("H100").Value = 2466.3333
$ws.Range("I100").Value = 2466.3333
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 2466.3333
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -1925.3333
$ws.Range("N100").ClearContents()
$ws.Range("H138").Value = 5399.8
$ws.Range("I138").Value = 5249.5
$ws.Range("J138").Value = 5500
$ws.Range("K138").Value = 15748.5
$ws.Range("L138").Value = 16500
$ws.Range("M138").Value = -10608.5
$ws.Range("N138").Value = -26780

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H58").Value = 20000
$ws.Range("J58").Value = 20000
$ws.Range("L58").Value = 20000
$ws.Range("N58").Value = -20860
$ws.Range("H74").Value = 3816.889
$ws.Range("I74").Value = 3816.889
$ws.Range("K74").Value = 3816.889
$ws.Range("M74").Value = -2942.889
$ws.Range("H77").Value = 3816.889
$ws.Range("I77").Value = 3816.889
$ws.Range("K77").Value = 19084.445
$ws.Range("M77").Value = -14716.445
$ws.Range("H110").Value = 699.75
$ws.Range("I110").Value = 699.75
$ws.Range("K110").Value = 699.75
$ws.Range("M110").Value = 1345.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").ClearContents()
$ws.Range("H112").Value = 25000
$ws.Range("J112").Value = 25000
$ws.Range("L112").Value = 25000
$ws.Range("N112").Value = -27954

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 943.0741
$ws.Range("I16").Value = 902.4231
$ws.Range("K16").Value = 902.4231
$ws.Range("M16").Value = -615.4231
$ws.Range("H21").Value = 40000
$ws.Range("I21").Value = 40000
$ws.Range("K21").Value = 40000
$ws.Range("M21").Value = -39765
$ws.Range("H42").Value = 23499.5
$ws.Range("I42").Value = 14000
$ws.Range("K42").Value = 14000
$ws.Range("M42").Value = -13407
$ws.Range("H99").Value = 731303.5
$ws.Range("I99").Value = 522050.3
$ws.Range("K99").Value = 522050.3
$ws.Range("M99").Value = -520552.3
$ws.Range("H113").Value = 943.0741
$ws.Range("I113").Value = 902.4231
$ws.Range("K113").Value = 902.4231
$ws.Range("M113").Value = 1267.5769
$ws.Range("H126").Value = 731303.5
$ws.Range("I126").Value = 522050.3
$ws.Range("K126").Value = 1566150.9
$ws.Range("M126").Value = -1563680.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 1100.2
$ws.Range("I114").Value = 1100.2
$ws.Range("J114").Value = 0
$ws.Range("K114").Value = 3300.6
$ws.Range("L114").Value = 0
$ws.Range("M114").Value = -46.60000000000036
$ws.Range("N114").ClearContents()
$ws.Range("H132").Value = 10000
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 90000
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -95060

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 36631
$ws.Range("I102").Value = 36631
$ws.Range("K102").Value = 36631
$ws.Range("M102").Value = -35009
$ws.Range("H126").Value = 11499.5
$ws.Range("I126").Value = 11499.5
$ws.Range("K126").Value = 34498.5
$ws.Range("M126").Value = -32028.5
$ws.Range("H132").Value = 1500
$ws.Range("I132").Value = 1000
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 3000
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -470
$ws.Range("N132").Value = -11060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2642.3333
$ws.Range("J22").Value = 2971
$ws.Range("L22").Value = 2971
$ws.Range("N22").Value = -3561
$ws.Range("H24").Value = 4000000
$ws.Range("J24").Value = 4000000
$ws.Range("L24").Value = 4000000
$ws.Range("N24").Value = -4000686
$ws.Range("H27").Value = 2642.3333
$ws.Range("J27").Value = 2971
$ws.Range("L27").Value = 2971
$ws.Range("N27").Value = -3185
$ws.Range("H61").Value = 950
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("H69").Value = 50000
$ws.Range("J69").Value = 50000
$ws.Range("L69").Value = 50000
$ws.Range("N69").Value = -51622
$ws.Range("H72").Value = 50000
$ws.Range("J72").Value = 50000
$ws.Range("L72").Value = 150000
$ws.Range("N72").Value = -158112
$ws.Range("H113").Value = 950
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("H136").Value = 2400
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 10799
$ws.Range("I32").Value = 11333.333
$ws.Range("J32").Value = 9997.5
$ws.Range("K32").Value = 11333.333
$ws.Range("L32").Value = 9997.5
$ws.Range("M32").Value = -11016.333
$ws.Range("N32").Value = -10631.5
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H80").Value = 22000
$ws.Range("I80").Value = 20000
$ws.Range("K80").Value = 20000
$ws.Range("M80").Value = -19002
$ws.Range("H83").Value = 22000
$ws.Range("I83").Value = 20000
$ws.Range("K83").Value = 60000
$ws.Range("M83").Value = -55008
$ws.Range("H126").Value = 2500
$ws.Range("I126").Value = 2500
$ws.Range("K126").Value = 7500
$ws.Range("M126").Value = -5030
$ws.Range("H132").Value = 2099
$ws.Range("I132").Value = 1918.8
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 5756.4
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -3226.4
$ws.Range("N132").Value = -14060
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H136").Value = 6447.625
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 6447.625
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 19342.875
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -24442.875
